$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (Volume/Number and week-range dates) ----
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# ---- Style reference cells (untouched row 14 covers styles 13/14/15) ----
$textStyleRef = $ws.Range("C14")   # style 13 (text placeholder style)
$numStyleRef  = $ws.Range("F14")   # style 14 (integer count style)
$pctStyleRef  = $ws.Range("L14")   # style 15 (percent-change style)

# ---- Data cell updates ----
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$textStyleRef.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$textStyleRef.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F15").Value2 = 2

$ws.Range("H15").Value2 = 0

$ws.Range("I15").Value2 = 2

$ws.Range("K15").Value2 = 0

$ws.Range("L15").Value2 = 100

$ws.Range("M15").Value2 = -33.333333333333

$ws.Range("N15").Value2 = -84.615384615384

$ws.Range("C16").Value2 = 4

$ws.Range("D16").Value2 = 1

$ws.Range("E16").Value2 = 300

$ws.Range("F16").Value2 = 16

$ws.Range("H16").Value2 = 100

$ws.Range("I16").Value2 = 18

$ws.Range("J16").Value2 = 9

$ws.Range("K16").Value2 = 100

$ws.Range("L16").Value2 = -25

$ws.Range("M16").Value2 = -28

$ws.Range("N16").Value2 = -83.636363636363

$ws.Range("C17").Value2 = 7

$ws.Range("D17").Value2 = 12

$ws.Range("E17").Value2 = -41.666666666666

$ws.Range("F17").Value2 = 19

$ws.Range("G17").Value2 = 27

$ws.Range("H17").Value2 = -29.629629629629

$ws.Range("I17").Value2 = 22

$ws.Range("J17").Value2 = 34

$ws.Range("K17").Value2 = -35.294117647058

$ws.Range("L17").Value2 = -38.888888888888

$ws.Range("M17").Value2 = 10

$ws.Range("N17").Value2 = -67.164179104477

$ws.Range("C18").Value2 = 7

$ws.Range("F18").Value2 = 15

$ws.Range("G18").Value2 = 5

$ws.Range("H18").Value2 = 200

$ws.Range("I18").Value2 = 16

$ws.Range("K18").Value2 = 220

$ws.Range("L18").Value2 = 23.076923076923

$ws.Range("M18").Value2 = -30.434782608695

$ws.Range("N18").Value2 = -87.692307692307

$ws.Range("D19").Value2 = 6

$ws.Range("E19").Value2 = -33.333333333333

$ws.Range("F19").Value2 = 19

$ws.Range("G19").Value2 = 26

$ws.Range("H19").Value2 = -26.923076923076

$ws.Range("I19").Value2 = 21

$ws.Range("J19").Value2 = 29

$ws.Range("K19").Value2 = -27.586206896551

$ws.Range("L19").Value2 = -51.162790697674

$ws.Range("M19").Value2 = -22.222222222222

$ws.Range("N19").Value2 = -58.823529411764

$ws.Range("D20").Value2 = 4

$ws.Range("E20").Value2 = -75

$ws.Range("G20").Value2 = 11

$ws.Range("H20").Value2 = -63.636363636363

$ws.Range("I20").Value2 = 4

$ws.Range("J20").Value2 = 12

$ws.Range("K20").Value2 = -66.666666666666

$ws.Range("L20").Value2 = -60

$ws.Range("M20").Value2 = -80

$ws.Range("N20").Value2 = -97.350993377483

$ws.Range("C21").Value2 = 23

$ws.Range("D21").Value2 = 23

$ws.Range("E21").Value2 = 0

$ws.Range("F21").Value2 = 76

$ws.Range("G21").Value2 = 79

$ws.Range("H21").Value2 = -3.79746835443

$ws.Range("I21").Value2 = 84

$ws.Range("J21").Value2 = 91

$ws.Range("K21").Value2 = -7.692307692307

$ws.Range("L21").Value2 = -34.375

$ws.Range("M21").Value2 = -28.813559322033

$ws.Range("N21").Value2 = -83.938814531548

$numStyleRef.Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C23").Value2 = 1

$ws.Range("E23").Value2 = 0

$ws.Range("F23").Value2 = 6

$ws.Range("H23").Value2 = 50

$ws.Range("I23").Value2 = 6

$ws.Range("J23").Value2 = 5

$ws.Range("K23").Value2 = 20

$ws.Range("M23").Value2 = 200

$ws.Range("C24").Value2 = 37

$ws.Range("D24").Value2 = 19

$ws.Range("E24").Value2 = 94.736842105263

$ws.Range("F24").Value2 = 111

$ws.Range("G24").Value2 = 99

$ws.Range("H24").Value2 = 12.121212121212

$ws.Range("I24").Value2 = 130

$ws.Range("J24").Value2 = 118

$ws.Range("K24").Value2 = 10.169491525423

$ws.Range("L24").Value2 = 15.04424778761

$ws.Range("M24").Value2 = 3.174603174603

$ws.Range("C25").Value2 = 25

$ws.Range("D25").Value2 = 7

$ws.Range("E25").Value2 = 257.142857142857

$ws.Range("F25").Value2 = 68

$ws.Range("G25").Value2 = 40

$ws.Range("H25").Value2 = 70

$ws.Range("I25").Value2 = 83

$ws.Range("J25").Value2 = 47

$ws.Range("K25").Value2 = 76.595744680851

$ws.Range("L25").Value2 = 72.916666666666

$ws.Range("C26").Value2 = 20

$ws.Range("D26").Value2 = 7

$ws.Range("E26").Value2 = 185.714285714286

$ws.Range("F26").Value2 = 55

$ws.Range("G26").Value2 = 40

$ws.Range("H26").Value2 = 37.5

$ws.Range("I26").Value2 = 69

$ws.Range("J26").Value2 = 49

$ws.Range("K26").Value2 = 40.816326530612

$ws.Range("L26").Value2 = 6.153846153846

$ws.Range("M26").Value2 = -30.30303030303

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$textStyleRef.Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D27").Value2 = 2

$ws.Range("E27").Value2 = -100

$ws.Range("F27").Value2 = 4

$ws.Range("G27").Value2 = 5

$ws.Range("H27").Value2 = -20

$ws.Range("I27").Value2 = 4

$ws.Range("J27").Value2 = 5

$ws.Range("K27").Value2 = -20

$ws.Range("L27").Value2 = 300

$ws.Range("D28").Value2 = 2

$ws.Range("E28").Value2 = 0

$ws.Range("G28").Value2 = 7

$ws.Range("H28").Value2 = -42.857142857142

$ws.Range("I28").Value2 = 5

$ws.Range("J28").Value2 = 8

$ws.Range("K28").Value2 = -37.5

$ws.Range("L28").Value2 = 0

$numStyleRef.Copy()
$ws.Range("D29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D29").Value2 = 2

$pctStyleRef.Copy()
$ws.Range("E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E29").Value2 = -100

$ws.Range("G29").Value2 = 2

$ws.Range("J29").Value2 = 3

$numStyleRef.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D30").Value2 = 2

$pctStyleRef.Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E30").Value2 = -100

$ws.Range("G30").Value2 = 2

$ws.Range("J30").Value2 = 3

$ws.Range("G31").Value2 = 2

$ws.Range("J31").Value2 = 2
